# Update the ranking table data (matrices scores recalculated, which shuffles
# the mat_rank-adjacent ordering for a couple of rows/people).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Colleen: only matrices score changes
$ws.Range("F2").Value = 13.4562512144795

# Row 3 - Annes: only matrices score changes
$ws.Range("F3").Value = 13.03225951465968

# Row 4 - now Bri (was Khushi)
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("D4").Value = "Bri"
$ws.Range("F4").Value = 8.49956431831203
$ws.Range("G4").Value = "Black or African American"

# Row 5 - Jewel: only matrices score changes
$ws.Range("F5").Value = 8.20087022370102

# Row 6 - now Khushi (was Bri)
$ws.Range("B6").Value = 22
$ws.Range("C6").Value = "608b14a312c099ac00b721b6"
$ws.Range("D6").Value = "Khushi"
$ws.Range("F6").Value = 8.064168822103696
$ws.Range("G6").Value = "Asian"

# Row 7 - Kellie: only matrices score changes
$ws.Range("F7").Value = 5.333931338090698

# Row 8 - Shadaisia: only matrices score changes
$ws.Range("F8").Value = 5.316202313826643

# Row 9 - Shaniek: only matrices score changes
$ws.Range("F9").Value = 5.294121455295787

# Row 10 - now Tina (was Lori)
$ws.Range("B10").Value = 34
$ws.Range("C10").Value = "5e96194b0a9fe909389e9f7b"
$ws.Range("D10").Value = "Tina"
$ws.Range("F10").Value = 4.302967855272668

# Row 11 - now Lori (was Tina)
$ws.Range("B11").Value = 35
$ws.Range("C11").Value = "6077db0613ce87b4a62a78f9"
$ws.Range("D11").Value = "Lori"
$ws.Range("F11").Value = 4.23075704731449

# Row 12 - Giana: only matrices score changes
$ws.Range("F12").Value = 2.016984074606204

# Row 13 - Nansi: only matrices score changes
$ws.Range("F13").Value = 1.054877676087834
